# Apply the edits described by the diff:
# 1) Title paragraph: font size 24 -> 28, remove the single underline.
# 2) Fix typo "FInish Line" -> "Finish Line".
# 3) Fix "decrease velocity of backward motion" -> "decrease velocity of forward motion"
#    only in the "Down Arrow" paragraph (leave the "Up Arrow" paragraph untouched).

$d = $word.ActiveDocument

# --- 1) Title formatting -------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleRange.Font.Size = 14
$titleRange.Font.SizeBi = 14
$titleRange.Font.Underline = 0

# --- 2) Typo fix: FInish Line -> Finish Line ------------------------------
$d.Content.Find.Execute("FInish Line", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Finish Line", 2)

# --- 3) Down Arrow paragraph: backward motion -> forward motion ----------
# Scope the Find to the "Down Arrow" paragraph only, after its leading tab,
# so the "Up Arrow" paragraph (identical trailing phrase) is left untouched.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Down Arrow: If stopped*") {
        $paraText = $p.Range.Text
        $tabIdx = $paraText.IndexOf([char]9)
        $bodyStart = $p.Range.Start + $tabIdx + 1
        $bodyRange = $d.Range($bodyStart, $p.Range.End)
        $bodyRange.Find.Execute("decrease velocity of backward motion", $true, $false, $false, $false, $false,
                                 $true, 1, $false, "decrease velocity of forward motion", 2)
    }
}
